$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("C1").Value = "email"
$ws.Range("A2").Value = "jonsmith"
$ws.Range("C2").Value = "NULL"
$ws.Range("B2").Value = "jon1234"

$ws.Range("D7").Select()
